$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/benefit-status"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
# Fixed Value for Extension.url mirrors the same URL used on the Metadata tab
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/benefit-status"
# Constraint(s) for the root "Extension" row is cleared
$elements.Range("AI2").Value = ""
